$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2:C2)
$ws.Range("A2").Value = " iAU_TC_ID_218"
$ws.Range("B2").Value = "@RegressionA Pre-Request Verify Elumina Login and Create Exam"
$ws.Range("C2").Value = "passed"

# Row 3 (A3:C3)
$ws.Range("A3").Value = " iAU_TC_ID_218"
$ws.Range("B3").Value = '@RegressionA Pre-Request "Validation of Delivery --> Add New Users"'
$ws.Range("C3").Value = "passed"

# Row 4 (A4:C4)
$ws.Range("A4").Value = "iAU_TC_ID_218"
$ws.Range("B4").Value = "@RegressionA Validation of Delivery--> Live Monitor - Candidate answer response Validation"
$ws.Range("C4").Value = "passed"

# Row 5 (A5:C5)
$ws.Range("A5").Value = "iAU_TC_ID_219A"
$ws.Range("B5").Value = "@RegressionA Validation of Delivery--> Live Monitor - Live Streaming page"
$ws.Range("C5").Value = "passed"
